$d = $word.ActiveDocument

# The paragraph "<id>p150v_1</id>" (the one referencing p150v_1, not the
# fig_p151r_1 one) is split across three runs with different formatting:
#   run1: "<id>"     (Courier New, color 7f6000, sz 18)
#   run2: "p150v_1"  (color 000000)
#   run3: "</id>"    (Courier New, color 7f6000, sz 18)
# Collapse them into a single run carrying run1's formatting, same as
# Word does when you Find/Replace text that spans multiple runs with
# identical replacement text - the match collapses onto the formatting
# of its first character.
$find = $d.Content.Find
$find.Execute("<id>p150v_1</id>", $false, $false, $false, $false, $false, `
               $true, 1, $false, "<id>p150v_1</id>", 2)
